$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric literal into a cell that lives in a column whose
# default/inherited format is Text ("@"), without the value getting coerced
# into a text shared-string. Toggling the format to General around the
# assignment keeps the stored value numeric; restoring "@" afterwards
# reproduces the original column styling (style index 2) on the new cell.
function Set-NumericInTextColumn($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = "@"
}

# --- Existing rows 5-7: status flips from RUNNING to DONE in column Q ---
$ws.Cells.Item(5, 17).Value = "DONE"
$ws.Cells.Item(6, 17).Value = "DONE"
$ws.Cells.Item(7, 17).Value = "DONE"

# --- New simulation rows 8-12 ---
$rows = @(
    @{ Row = 8;  Coeff = "0.3";  Folder = "param01seed0.30"; Status = "DONE" },
    @{ Row = 9;  Coeff = "0.35"; Folder = "param01seed0.35"; Status = "DONE" },
    @{ Row = 10; Coeff = "0.4";  Folder = "param01seed0.40"; Status = "DONE" },
    @{ Row = 11; Coeff = "0.45"; Folder = "param01seed0.45"; Status = "RUNNING" },
    @{ Row = 12; Coeff = "0.5";  Folder = "param01seed0.50"; Status = "RUNNING" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = "BA"
    $ws.Cells.Item($row, 2).Value = "n=100000 m=2"
    $ws.Cells.Item($row, 3).Value = "m=100000 k=4 p=0.7"
    Set-NumericInTextColumn $row 4 999999

    $ws.Cells.Item($row, 5).Formula = "=" + $r.Coeff + "*(F" + $row + "+G" + $row + "+H" + $row + "+I" + $row + ")"

    Set-NumericInTextColumn $row 6 0.45
    Set-NumericInTextColumn $row 7 0.09
    Set-NumericInTextColumn $row 8 0.0225
    Set-NumericInTextColumn $row 9 0

    $ws.Cells.Item($row, 10).Value = "52/6"
    $ws.Cells.Item($row, 11).Value = "365/77"
    $ws.Cells.Item($row, 12).Value = "365/77"
    $ws.Cells.Item($row, 13).Value = "52/12"
    $ws.Cells.Item($row, 14).Value = "12/25"
    $ws.Cells.Item($row, 15).Value = "12/25"
    $ws.Cells.Item($row, 16).Value = $r.Folder
    $ws.Cells.Item($row, 17).Value = $r.Status
}
